# ---------------------------------------------------------------------------
# Tactile Tabletop Data - Level 3 Character Cards: content + layout refresh
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 1-5: values only shuffle columns / drop "1 Skillpoint"; write row-major ---
$row1 = New-Object 'object[,]' 1,16
$row1[0,0] = "Top Ability Name"
$row1[0,1] = "Top Ability Target"
$row1[0,2] = "Top Ability Duration"
$row1[0,3] = "Top Ability Die Roll/Scaler"
$row1[0,4] = "Top Ability Rules"
$row1[0,5] = "Top Ability Following Card Action"
$row1[0,6] = "Bottom Ability Name"
$row1[0,7] = "Bottom Ability Target"
$row1[0,8] = "Bottom Ability Duration"
$row1[0,9] = "Bottom Ability Die Roll/Scaler"
$row1[0,10] = "Bottom Ability Rules"
$row1[0,11] = "Bottom Ability Following Card Action"
$row1[0,12] = "Passives"
$row1[0,13] = "Requirements"
$row1[0,14] = $null
$row1[0,15] = "General Notes"
$ws.Range("A1:P1").Value = $row1

$row2 = New-Object 'object[,]' 1,14
$row2[0,0] = "Upgrade"
$row2[0,1] = "Self"
$row2[0,2] = "X Rnds"
$row2[0,3] = "X = Level   Y = Influence"
$row2[0,4] = "Exhaust a card. Add Y to each Defense or Attack roll you make"
$row2[0,5] = "Discard"
$row2[0,6] = "Built from Scraps"
$row2[0,7] = "Self"
$row2[0,8] = "Instant"
$row2[0,9] = "X = Exhausted cards from Discard"
$row2[0,10] = "Spend X to perform the following: `n1) Add 25 feet to your movement`n3) Attack target enemy, they cannot defend`n5) Take an additional turn after this one"
$row2[0,11] = "Exhaust"
$row2[0,12] = $null
$row2[0,13] = "3 Craftsmanship, 2 Knowledge"
$ws.Range("A2:N2").Value = $row2

$row3 = New-Object 'object[,]' 1,14
$row3[0,0] = "Well of Life"
$row3[0,1] = "Self"
$row3[0,2] = "1 Rnd"
$row3[0,3] = "X = Influence"
$row3[0,4] = "Action can only be used by itself. At start of next turn, recover all discarded cards, Heal for X, and add X to defence die for 1 rnd"
$row3[0,5] = "Exhaust"
$row3[0,6] = "Guided Strike"
$row3[0,7] = "Enemy"
$row3[0,8] = "1 Rnd"
$row3[0,9] = "Attack"
$row3[0,10] = "Attack target once this turn, again at start of next turn"
$row3[0,11] = "Discard"
$row3[0,12] = $null
$row3[0,13] = "2 Spirituality, 3 Charisma"
$ws.Range("A3:N3").Value = $row3

$row4 = New-Object 'object[,]' 1,14
$row4[0,0] = "Neck Cracker"
$row4[0,1] = "Enemy"
$row4[0,2] = "Instant"
$row4[0,3] = "Attack    X = Level"
$row4[0,4] = "Make an Unarmed attack, add X to the Attack Value. If you deal damage Attack Target Enemy. Target Enemy cannot use movement on their next turn."
$row4[0,5] = "Discard"
$row4[0,6] = "Shatter Weapon"
$row4[0,7] = "Enemy"
$row4[0,8] = "Instant"
$row4[0,9] = "Attack    X = Level"
$row4[0,10] = "Make an Unarmed attack against an enemy and add X to the Attack Value.Lose X life. If you deal damage, you may choose to destroy one of their equipped weapons"
$row4[0,11] = "Exhaust"
$row4[0,12] = $null
$row4[0,13] = "4 STR and 2 VIG"
$ws.Range("A4:N4").Value = $row4

$row5 = New-Object 'object[,]' 1,26
$row5[0,0] = "Lightning Strike"
$row5[0,1] = "Enemy"
$row5[0,2] = "Instant"
$row5[0,3] = "Attack    X = Level"
$row5[0,4] = "Attack Target Enemy twice and add X to each Attack Value. Discard a card."
$row5[0,5] = "Hand"
$row5[0,6] = "Reliable Strike"
$row5[0,7] = "Enemy"
$row5[0,8] = "Instant"
$row5[0,9] = "Attack    X = Level"
$row5[0,10] = "Attack target, add X to your Attack Value. If no damage is dealt then return this card to your hand."
$row5[0,11] = "Discard"
$row5[0,12] = $null
$row5[0,13] = "Quick Strike"
$row5[0,14] = "Enemy"
$row5[0,15] = "Instant"
$row5[0,16] = "Attack"
$row5[0,17] = "Discard a card"
$row5[0,18] = "Hand"
$row5[0,19] = "Basic Strike"
$row5[0,20] = "Enemy"
$row5[0,21] = "Instant"
$row5[0,22] = "Attack"
$row5[0,23] = "Attack target"
$row5[0,24] = "Discard"
$row5[0,25] = "1 Health"
$ws.Range("A5:Z5").Value = $row5

# --- Rows 6-11: brand-new card rows. New text is entered cell-by-cell in the
#     exact order the original author typed it, so the workbook's shared-string
#     table is appended to in the same sequence as the source file.
$ws.Range("A6").Value = "Turn Warp"
$ws.Range("D6").Value = "X = Level"
$ws.Range("E6").Value = "Lose X life. Take an extra turn after this turn ends and then lose X life again. After taking your extra turn, you must skip your next turn. "
$ws.Range("G6").Value = "Extra Steps"
$ws.Range("N6").Value = "6 INT"
$ws.Range("K6").Value = "Perform X non-action abilities, and lose 1 life for each one. "
$ws.Range("A7").Value = "Organize"
$ws.Range("E7").Value = "Decide, layout, and expose the actions for your next turn. Everyone/everything around you knows your next turn. Recover X life, return X cards from discard, and add X to your next roll."
$ws.Range("G7").Value = "Scheme"
$ws.Range("K7").Value = "During the next round, you may play either of your two actions at any point not during an active turn. At the end of that round, play any actions you have not already. "
$ws.Range("L7").Value = "Return to Hand"
$ws.Range("A8").Value = "Richochet"
$ws.Range("G8").Value = "Near Miss"
$ws.Range("K8").Value = "If your attack value does not damage the targeted enemy/enemies, they take X damage."
$ws.Range("E8").Value = "If your next attack value does not damage a targetted enemy, add X to the attak value and target another enemy in range. Repeat this proces until damage is dealt, or there are no new enemies you can target. "
$ws.Range("A9").Value = "Gravity Well"
$ws.Range("B9").Value = "Area"
$ws.Range("D9").Value = "X = Level "
$ws.Range("G9").Value = "Inversion"
$ws.Range("E9").Value = "Target area 5 feet by 5 feet. All Enemies in and adjacanet to that area cannot move more than 5 feet from targetted area unless they take and suceed an influence check against you."
$ws.Range("K9").Value = "Target area 5 feet by 5 feet. Enemies cannot move within 5 feet of the targetted area unless they make and suceed an influence roll against you."
$ws.Range("A10").Value = "Hookshot"
$ws.Range("E10").Value = "Target area within 45 feet and jump to that location"
$ws.Range("G10").Value = "Lasso"
$ws.Range("H10").Value = "Ally"
$ws.Range("K10").Value = "Target Ally within 45 feet and pull them to you."
$ws.Range("A11").Value = " Touch"

# --- Remaining rows 6-10 cells that only reuse already-existing shared strings ---
$ws.Range("B6").Value = "Self"
$ws.Range("C6").Value = "Instant"
$ws.Range("F6").Value = "Exhaust"
$ws.Range("H6").Value = "Self"
$ws.Range("I6").Value = "Instant"
$ws.Range("L6").Value = "Exhaust"
$ws.Range("B7").Value = "Self"
$ws.Range("C7").Value = "Instant"
$ws.Range("F7").Value = "Exhaust"
$ws.Range("H7").Value = "Self"
$ws.Range("I7").Value = "Instant"
$ws.Range("B8").Value = "Enemy"
$ws.Range("C8").Value = "X Rnds"
$ws.Range("F8").Value = "Discard"
$ws.Range("H8").Value = "Self"
$ws.Range("I8").Value = "X Rnds"
$ws.Range("L8").Value = "Discard"
$ws.Range("C9").Value = "X Rnds"
$ws.Range("F9").Value = "Discard"
$ws.Range("I9").Value = "X Rnds"
$ws.Range("L9").Value = "Discard"
$ws.Range("C10").Value = "Instant"
$ws.Range("F10").Value = "Discard"
$ws.Range("I10").Value = "Instant"
$ws.Range("L10").Value = "Discard"

# --- Column widths (A:N) matching the new, wider layout ---
$ws.Columns.Item(1).ColumnWidth = 16.333333333333332
$ws.Columns.Item(2).ColumnWidth = 17.5
$ws.Columns.Item(3).ColumnWidth = 18.666666666666668
$ws.Columns.Item(4).ColumnWidth = 24.833333333333332
$ws.Columns.Item(5).ColumnWidth = 19.166666666666668
$ws.Columns.Item(6).ColumnWidth = 30.333333333333332
$ws.Columns.Item(7).ColumnWidth = 20.666666666666668
$ws.Columns.Item(8).ColumnWidth = 21.666666666666668
$ws.Columns.Item(9).ColumnWidth = 21.333333333333332
$ws.Columns.Item(10).ColumnWidth = 28.0
$ws.Columns.Item(11).ColumnWidth = 20.833333333333332
$ws.Columns.Item(12).ColumnWidth = 34.333333333333336
$ws.Columns.Item(13).ColumnWidth = 12.833333333333334
$ws.Columns.Item(14).ColumnWidth = 32.0

# --- Active selection lands on the newly-added A11 cell ---
$ws.Range("A11").Select()
